$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per row, exactly matching the target diff.
# For Price (column D) values that look like plain numbers, force text
# formatting first so Excel does not silently coerce them to numeric values
# (the sheet stores these as text, e.g. "1.000", "29.282.69", etc.).

# Row 2
$ws.Range("D2").Value = "29.282.69"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "1.870.54"
$ws.Range("E3").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7086"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.60"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07787"
$ws.Range("E8").Value = "  +1.60%  "

# Row 9
$ws.Range("E9").Value = "  -0.70%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.97"
$ws.Range("E10").Value = "  +1.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08398"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12
$ws.Range("D12").Value = "1.874.44"
$ws.Range("E12").Value = "  +0.41%  "

# Row 13
$ws.Range("E13").Value = "  +0.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7098"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("E15").Value = "  -0.41%  "

# Row 16
$ws.Range("D16").Value = "29.294.99"
$ws.Range("E16").Value = "  +0.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.059"
$ws.Range("E17").Value = "  +1.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008186"
$ws.Range("E18").Value = "  +4.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.47"
$ws.Range("E19").Value = "  -1.77%  "

# Row 20
$ws.Range("E20").Value = "  +0.93%  "

# Row 21
$ws.Range("D21").Value = "2.114.57"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.745"
$ws.Range("E23").Value = "  -1.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.0000"
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("E25").Value = "  -3.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.22"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.002"
$ws.Range("E27").Value = "  +0.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.42"
$ws.Range("E28").Value = "  -0.32%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.387"
$ws.Range("E30").Value = "  -0.34%  "

# Row 31
$ws.Range("E31").Value = "  -1.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.284"
$ws.Range("E32").Value = "  +0.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05346"
$ws.Range("E33").Value = "  +3.55%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.937"
$ws.Range("E34").Value = "  +1.20%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.175"
$ws.Range("E35").Value = "  +0.71%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7453"
$ws.Range("E36").Value = "  -6.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.698"
$ws.Range("E37").Value = "  +0.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01868"
$ws.Range("E38").Value = "  +0.41%  "

# Row 39
$ws.Range("D39").Value = "1.229.69"
$ws.Range("E39").Value = "  +5.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.726"
$ws.Range("E40").Value = "  +0.62%  "

# Row 41
$ws.Range("E41").Value = "  +3.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8841"
$ws.Range("E42").Value = "  -1.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.87"
$ws.Range("E43").Value = "  +5.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.25"
$ws.Range("E44").Value = "  -1.85%  "

# Row 45
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("D46").Value = "2.010.98"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5190"
$ws.Range("E47").Value = "  +0.25%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000122"
$ws.Range("E48").Value = "  +2.49%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.787"
$ws.Range("E49").Value = "  +0.52%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.417"
$ws.Range("E50").Value = "  +0.74%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4305"
$ws.Range("E51").Value = "  +0.21%  "
